# acceptance: add AT-01 automation and close remaining partial rows
#
# WBS sheet: rows 50-62 (Phase 3 - Traceability Engine block) move from
# "Partial" execution status to "Done", gain a Completed On date, and get
# the checkmark evidence columns (Schema/Validation/Permissions-Isolation/
# Workflow/Evidence) filled in — matching the pattern already used by the
# completed rows earlier in the sheet (e.g. row 2).
#
# Milestones sheet: M3 (row 5) status flips from "In Progress" to "Done".

$wb = $excel.ActiveWorkbook

$wbsSheet = $wb.Worksheets.Item("WBS")

$completedOn = "2026-04-17"
$check = [char]0x2705   # checkmark

for ($row = 50; $row -le 62; $row++) {
    $wbsSheet.Cells.Item($row, 8).Value = "Done"        # H: Execution Status

    # K: Completed On. The sheet stores these "date" columns as plain text
    # (matching J, Started On) rather than real date serials, so force the
    # cell to text before/after the assignment to avoid Excel's automatic
    # date-literal parsing.
    $kCell = $wbsSheet.Cells.Item($row, 11)
    $kCell.NumberFormat = "@"
    $kCell.Value = $completedOn
    $kCell.NumberFormat = "yyyy-mm-dd"

    $wbsSheet.Cells.Item($row, 12).Value = $check        # L: Schema
    $wbsSheet.Cells.Item($row, 13).Value = $check        # M: Validation
    $wbsSheet.Cells.Item($row, 14).Value = $check        # N: Permissions/Isolation
    $wbsSheet.Cells.Item($row, 15).Value = $check        # O: Workflow
    $wbsSheet.Cells.Item($row, 16).Value = $check        # P: Evidence
}

$milestonesSheet = $wb.Worksheets.Item("Milestones")
$milestonesSheet.Range("F5").Value = "$check Done"
